# Apply updated plotting-script values to each Fold sheet (rows 2-4, cols B-O)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$arr = New-Object "object[,]" 3,14
$arr[0,0]=5.45195; $arr[0,1]=5.89945; $arr[0,2]=37.50169999999999; $arr[0,3]=0.17115; $arr[0,4]=2219.401700000001; $arr[0,5]=2014.1057; $arr[0,6]=205.29605; $arr[0,7]=205.29605; $arr[0,8]=0; $arr[0,9]=2798.105; $arr[0,10]=2014.125; $arr[0,11]=783.98; $arr[0,12]=345.6012; $arr[0,13]=438.3792
$arr[1,0]=5.789; $arr[1,1]=9.572; $arr[1,2]=28.079; $arr[1,3]=0; $arr[1,4]=2327.595; $arr[1,5]=2122.899; $arr[1,6]=204.696; $arr[1,7]=193.677; $arr[1,8]=11.019; $arr[1,9]=2476.9342; $arr[1,10]=2122.95; $arr[1,11]=353.9842; $arr[1,12]=308.7674; $arr[1,13]=45.2168
$arr[2,0]=6.746; $arr[2,1]=8.249; $arr[2,2]=41.042; $arr[2,3]=0; $arr[2,4]=2510.203; $arr[2,5]=2280.003; $arr[2,6]=230.2; $arr[2,7]=64.409; $arr[2,8]=165.791; $arr[2,9]=2467.796; $arr[2,10]=2280.003; $arr[2,11]=187.793; $arr[2,12]=176.6258; $arr[2,13]=11.1674
$ws.Range("B2:O4").Value = $arr

$ws = $wb.Worksheets.Item("Fold_2")
$arr = New-Object "object[,]" 3,14
$arr[0,0]=5.591900000000001; $arr[0,1]=6.396350000000001; $arr[0,2]=36.86314999999999; $arr[0,3]=0.1292; $arr[0,4]=2253.7439; $arr[0,5]=2044.3985; $arr[0,6]=209.34555; $arr[0,7]=209.34555; $arr[0,8]=0; $arr[0,9]=2285.459; $arr[0,10]=2044.411; $arr[0,11]=241.048; $arr[0,12]=241.048; $arr[0,13]=0
$arr[1,0]=5.963; $arr[1,1]=9.461; $arr[1,2]=28.427; $arr[1,3]=0; $arr[1,4]=2359.512; $arr[1,5]=2141.396; $arr[1,6]=218.116; $arr[1,7]=201.182; $arr[1,8]=16.934; $arr[1,9]=2347.4208; $arr[1,10]=2141.429; $arr[1,11]=205.9918; $arr[1,12]=205.9918; $arr[1,13]=0
$arr[2,0]=6.7; $arr[2,1]=8.094; $arr[2,2]=47.574; $arr[2,3]=0; $arr[2,4]=2580.234; $arr[2,5]=2317.715; $arr[2,6]=262.519; $arr[2,7]=54.052; $arr[2,8]=208.467; $arr[2,9]=2395.058; $arr[2,10]=2317.715; $arr[2,11]=77.343; $arr[2,12]=77.343; $arr[2,13]=0
$ws.Range("B2:O4").Value = $arr

$ws = $wb.Worksheets.Item("Fold_3")
$arr = New-Object "object[,]" 3,14
$arr[0,0]=5.43975; $arr[0,1]=6.10355; $arr[0,2]=37.0707; $arr[0,3]=0.17115; $arr[0,4]=2222.8903; $arr[0,5]=2017.33845; $arr[0,6]=205.55185; $arr[0,7]=205.55185; $arr[0,8]=0; $arr[0,9]=2804.228; $arr[0,10]=2017.357; $arr[0,11]=786.8710000000001; $arr[0,12]=359.913; $arr[0,13]=426.9582
$arr[1,0]=5.483; $arr[1,1]=7.719; $arr[1,2]=32.87; $arr[1,3]=0; $arr[1,4]=2320.004; $arr[1,5]=2051.97; $arr[1,6]=268.033; $arr[1,7]=220.571; $arr[1,8]=47.462; $arr[1,9]=2618.863; $arr[1,10]=2051.924; $arr[1,11]=566.939; $arr[1,12]=353.6808; $arr[1,13]=213.2584
$arr[2,0]=6.231; $arr[2,1]=8.112; $arr[2,2]=43.556; $arr[2,3]=0; $arr[2,4]=2573.291; $arr[2,5]=2233.933; $arr[2,6]=339.358; $arr[2,7]=119.435; $arr[2,8]=219.923; $arr[2,9]=2490.067; $arr[2,10]=2233.933; $arr[2,11]=256.134; $arr[2,12]=190.343; $arr[2,13]=65.791
$ws.Range("B2:O4").Value = $arr

$ws = $wb.Worksheets.Item("Fold_4")
$arr = New-Object "object[,]" 3,14
$arr[0,0]=5.734050000000002; $arr[0,1]=6.146849999999999; $arr[0,2]=39.35619999999999; $arr[0,3]=0.04355; $arr[0,4]=2295.672100000001; $arr[0,5]=2069.291249999999; $arr[0,6]=226.3809; $arr[0,7]=226.3809; $arr[0,8]=0; $arr[0,9]=2185.4816; $arr[0,10]=2069.311; $arr[0,11]=116.1706; $arr[0,12]=80.70360000000001; $arr[0,13]=35.467
$arr[1,0]=6.113; $arr[1,1]=9.366; $arr[1,2]=29.606; $arr[1,3]=0; $arr[1,4]=2390.136; $arr[1,5]=2164.023; $arr[1,6]=226.112; $arr[1,7]=213.845; $arr[1,8]=12.267; $arr[1,9]=2231.9788; $arr[1,10]=2163.984; $arr[1,11]=67.9948; $arr[1,12]=67.9948; $arr[1,13]=0
$arr[2,0]=6.7; $arr[2,1]=8.094; $arr[2,2]=47.574; $arr[2,3]=0; $arr[2,4]=2580.234; $arr[2,5]=2317.715; $arr[2,6]=262.519; $arr[2,7]=54.052; $arr[2,8]=208.467; $arr[2,9]=2328.9178; $arr[2,10]=2317.715; $arr[2,11]=11.2028; $arr[2,12]=11.2028; $arr[2,13]=0
$ws.Range("B2:O4").Value = $arr

$ws = $wb.Worksheets.Item("Fold_5")
$arr = New-Object "object[,]" 3,14
$arr[0,0]=5.656750000000001; $arr[0,1]=5.934599999999999; $arr[0,2]=39.77945; $arr[0,3]=0.16955; $arr[0,4]=2266.771; $arr[0,5]=2056.2701; $arr[0,6]=210.50085; $arr[0,7]=210.50085; $arr[0,8]=0; $arr[0,9]=2393.0558; $arr[0,10]=2056.285; $arr[0,11]=336.7708; $arr[0,12]=152.2866; $arr[0,13]=184.4844
$arr[1,0]=6.099; $arr[1,1]=7.875; $arr[1,2]=33.748; $arr[1,3]=0; $arr[1,4]=2371.432; $arr[1,5]=2135.853; $arr[1,6]=235.579; $arr[1,7]=199.385; $arr[1,8]=36.195; $arr[1,9]=2332.4364; $arr[1,10]=2135.884; $arr[1,11]=196.5524; $arr[1,12]=126.9988; $arr[1,13]=69.5534
$arr[2,0]=7.162; $arr[2,1]=7.283; $arr[2,2]=49.794; $arr[2,3]=0; $arr[2,4]=2576.937; $arr[2,5]=2356.575; $arr[2,6]=220.362; $arr[2,7]=121.688; $arr[2,8]=98.675; $arr[2,9]=2467.0234; $arr[2,10]=2356.575; $arr[2,11]=110.4484; $arr[2,12]=44.6264; $arr[2,13]=65.8222
$ws.Range("B2:O4").Value = $arr
